$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 44140
$ws.Range("B4").Value = "Limette"
$ws.Range("B5").Value = 4

$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("B8").Value = 90
$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("B9").Value = 110
$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("B10").Value = 110
$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("B11").Value = 90

$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""

$ws.Range("A19").Value = "Study Report"
$ws.Range("B19").Value = "High Fidelity Prototype"

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B20").Select()
